# Adds a new "2022-Q4" worksheet (inserted right after the summary sheet
# "总计" and before "2021-Q3"), fills it with the quarter's fund-holding
# data, and updates the "总计" summary sheet with a new top data row for
# 2022-Q4 (existing rows shift down by one).
#
# The other quarter sheets (2021-Q3 / 2021-Q2 / 2021-Q1 / 2020-Q4) keep
# their original data; only their tab position shifts right by one to make
# room for the newly inserted "2022-Q4" sheet.

$wb = $excel.ActiveWorkbook

function Style-Header($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1  # xlEdgeLeft
    $cell.Borders.Item(8).LineStyle = 1  # xlEdgeTop
    $cell.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
    $cell.Borders.Item(10).LineStyle = 1 # xlEdgeRight
}

function Style-Index($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1. Insert the brand-new "2022-Q4" worksheet before "2021-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# Header row (note: this sheet's amount column header is "基金规模",
# unlike the other quarter sheets which use "基金金额").
Set-TextValue $newSheet.Cells.Item(1,2) "基金代码"
Set-TextValue $newSheet.Cells.Item(1,3) "基金名称"
Set-TextValue $newSheet.Cells.Item(1,4) "基金规模"
Set-TextValue $newSheet.Cells.Item(1,5) "股票总仓位"
Set-TextValue $newSheet.Cells.Item(1,6) "仓位占比"
Set-TextValue $newSheet.Cells.Item(1,7) "持有市值(亿元)"
Set-TextValue $newSheet.Cells.Item(1,8) "仓位排名"
foreach ($col in 2..8) {
    Style-Header $newSheet.Cells.Item(1,$col)
}

$fundRows = @(
    @(0, "003956", "南方产业智选股票",                  "3.60", "85.80", "4.89", "0.1760", 5),
    @(1, "003413", "华泰柏瑞新经济沪港深混合",             "1.44", "94.26", "7.86", "0.1132", 2),
    @(2, "671010", "西部利得策略优选混合A",                "1.88", "92.90", "5.35", "0.1006", 9),
    @(3, "010204", "中银港股通优势成长股票",               "2.98", "86.92", "2.35", "0.0700", 10),
    @(4, "011355", "华泰柏瑞港股通时代机遇混合A",          "0.70", "94.61", "7.17", "0.0502", 2),
    @(5, "015143", "中欧智能制造混合A",                   "1.54", "75.37", "2.52", "0.0388", 9),
    @(6, "460010", "华泰柏瑞亚洲领导企业混合（QDII）",     "0.52", "97.17", "6.07", "0.0316", 3),
    @(7, "011356", "华泰柏瑞港股通时代机遇混合C",          "0.39", "94.61", "7.17", "0.0280", 2),
    @(8, "011060", "西部利得策略优选混合C",                "0.48", "92.90", "5.35", "0.0257", 9),
    @(9, "015144", "中欧智能制造混合C",                   "0.67", "75.37", "2.52", "0.0169", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $idxCell = $newSheet.Cells.Item($r,1)
    $idxCell.Value = $row[0]
    Style-Index $idxCell

    Set-TextValue $newSheet.Cells.Item($r,2) $row[1]
    Set-TextValue $newSheet.Cells.Item($r,3) $row[2]
    Set-TextValue $newSheet.Cells.Item($r,4) $row[3]
    Set-TextValue $newSheet.Cells.Item($r,5) $row[4]
    Set-TextValue $newSheet.Cells.Item($r,6) $row[5]
    Set-TextValue $newSheet.Cells.Item($r,7) $row[6]
    $newSheet.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q4 right
#    after the header, push the existing quarters down, and keep the
#    running index (column A) sequential.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()

$summary.Cells.Item(2,1).Value = 0
Style-Index $summary.Cells.Item(2,1)
Set-TextValue $summary.Cells.Item(2,2) "2022-Q4"
$summary.Cells.Item(2,3).Value = 10
$summary.Cells.Item(2,4).Value = 0.65

# Renumber the running index for the rows that shifted down.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
